$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

# Insert a new column before column F (shifts F.. to G.. )
$ws.Columns("F:F").Insert()

# Populate the newly inserted column F
$ws.Cells.Item(1, 6).Value = "type_part"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 6).Value = "WALL_AS7"
}

# Update selection to match the final state (active cell F7)
$ws.Range("F7").Select() | Out-Null
